# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns land right after the existing last column (AC), i.e. AD:AF,
# growing the used range from A1:AC51 to A1:AF51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the three new header cells (AD1:AF1) the same look as the existing
# header row (bold font, thin border, centered/top aligned) by cloning the
# format from the adjacent header cell (AC1) instead of building a new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-51, including the stray repeated-header row 51) gets the
# same team record: 88 wins, 74 losses, 0 ties.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
